$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$cell = $t.Cell(1, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("74+25=", $true, $false, $false, $false, $false, $true, 1, $false, "41-21=", 2) | Out-Null
$cell = $t.Cell(1, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("21+34=", $true, $false, $false, $false, $false, $true, 1, $false, "4+59=", 2) | Out-Null
$cell = $t.Cell(1, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("30-10=", $true, $false, $false, $false, $false, $true, 1, $false, "85-62=", 2) | Out-Null
$cell = $t.Cell(1, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("49+31=", $true, $false, $false, $false, $false, $true, 1, $false, "80+17=", 2) | Out-Null
$cell = $t.Cell(1, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("99-23=", $true, $false, $false, $false, $false, $true, 1, $false, "24+38=", 2) | Out-Null
$cell = $t.Cell(2, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("90-82=", $true, $false, $false, $false, $false, $true, 1, $false, "52+31=", 2) | Out-Null
$cell = $t.Cell(2, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("64+23=", $true, $false, $false, $false, $false, $true, 1, $false, "61+13=", 2) | Out-Null
$cell = $t.Cell(2, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("22+21=", $true, $false, $false, $false, $false, $true, 1, $false, "8+64=", 2) | Out-Null
$cell = $t.Cell(2, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("72+15=", $true, $false, $false, $false, $false, $true, 1, $false, "48+6=", 2) | Out-Null
$cell = $t.Cell(2, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("57-45=", $true, $false, $false, $false, $false, $true, 1, $false, "42+28=", 2) | Out-Null
$cell = $t.Cell(3, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("1+23=", $true, $false, $false, $false, $false, $true, 1, $false, "35-9=", 2) | Out-Null
$cell = $t.Cell(3, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("53+22=", $true, $false, $false, $false, $false, $true, 1, $false, "86-17=", 2) | Out-Null
$cell = $t.Cell(3, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("6+67=", $true, $false, $false, $false, $false, $true, 1, $false, "26+53=", 2) | Out-Null
$cell = $t.Cell(3, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("55+21=", $true, $false, $false, $false, $false, $true, 1, $false, "44-32=", 2) | Out-Null
$cell = $t.Cell(3, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("4+42=", $true, $false, $false, $false, $false, $true, 1, $false, "89-36=", 2) | Out-Null
$cell = $t.Cell(4, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("6+18=", $true, $false, $false, $false, $false, $true, 1, $false, "11+13=", 2) | Out-Null
$cell = $t.Cell(4, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("11+28=", $true, $false, $false, $false, $false, $true, 1, $false, "34+39=", 2) | Out-Null
$cell = $t.Cell(4, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("32+58=", $true, $false, $false, $false, $false, $true, 1, $false, "95-77=", 2) | Out-Null
$cell = $t.Cell(4, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("68-4=", $true, $false, $false, $false, $false, $true, 1, $false, "91-37=", 2) | Out-Null
$cell = $t.Cell(4, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("72-40=", $true, $false, $false, $false, $false, $true, 1, $false, "86-22=", 2) | Out-Null
$cell = $t.Cell(5, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("45+14=", $true, $false, $false, $false, $false, $true, 1, $false, "14+42=", 2) | Out-Null
$cell = $t.Cell(5, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("34+12=", $true, $false, $false, $false, $false, $true, 1, $false, "75+24=", 2) | Out-Null
$cell = $t.Cell(5, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("43-2=", $true, $false, $false, $false, $false, $true, 1, $false, "72-71=", 2) | Out-Null
$cell = $t.Cell(5, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("18+79=", $true, $false, $false, $false, $false, $true, 1, $false, "68-23=", 2) | Out-Null
$cell = $t.Cell(5, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("9+57=", $true, $false, $false, $false, $false, $true, 1, $false, "29-16=", 2) | Out-Null
$cell = $t.Cell(6, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("33-5=", $true, $false, $false, $false, $false, $true, 1, $false, "87-76=", 2) | Out-Null
$cell = $t.Cell(6, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("40+40=", $true, $false, $false, $false, $false, $true, 1, $false, "40+46=", 2) | Out-Null
$cell = $t.Cell(6, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("94-31=", $true, $false, $false, $false, $false, $true, 1, $false, "50-46=", 2) | Out-Null
$cell = $t.Cell(6, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("11+61=", $true, $false, $false, $false, $false, $true, 1, $false, "79-42=", 2) | Out-Null
$cell = $t.Cell(6, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("61+22=", $true, $false, $false, $false, $false, $true, 1, $false, "66-33=", 2) | Out-Null
$cell = $t.Cell(7, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("72+20=", $true, $false, $false, $false, $false, $true, 1, $false, "86-19=", 2) | Out-Null
$cell = $t.Cell(7, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("76-72=", $true, $false, $false, $false, $false, $true, 1, $false, "90-64=", 2) | Out-Null
$cell = $t.Cell(7, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("54+23=", $true, $false, $false, $false, $false, $true, 1, $false, "63+20=", 2) | Out-Null
$cell = $t.Cell(7, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("26+73=", $true, $false, $false, $false, $false, $true, 1, $false, "15-0=", 2) | Out-Null
$cell = $t.Cell(7, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("25+59=", $true, $false, $false, $false, $false, $true, 1, $false, "62+13=", 2) | Out-Null
$cell = $t.Cell(8, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("1+43=", $true, $false, $false, $false, $false, $true, 1, $false, "12+34=", 2) | Out-Null
$cell = $t.Cell(8, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("94-14=", $true, $false, $false, $false, $false, $true, 1, $false, "9+19=", 2) | Out-Null
$cell = $t.Cell(8, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("16+55=", $true, $false, $false, $false, $false, $true, 1, $false, "67-26=", 2) | Out-Null
$cell = $t.Cell(8, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("65+24=", $true, $false, $false, $false, $false, $true, 1, $false, "32+40=", 2) | Out-Null
$cell = $t.Cell(8, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("81-74=", $true, $false, $false, $false, $false, $true, 1, $false, "57+4=", 2) | Out-Null
$cell = $t.Cell(9, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("79-36=", $true, $false, $false, $false, $false, $true, 1, $false, "20+46=", 2) | Out-Null
$cell = $t.Cell(9, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("38-14=", $true, $false, $false, $false, $false, $true, 1, $false, "27-16=", 2) | Out-Null
$cell = $t.Cell(9, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("47+50=", $true, $false, $false, $false, $false, $true, 1, $false, "84-47=", 2) | Out-Null
$cell = $t.Cell(9, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("96-36=", $true, $false, $false, $false, $false, $true, 1, $false, "16-8=", 2) | Out-Null
$cell = $t.Cell(9, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("35+29=", $true, $false, $false, $false, $false, $true, 1, $false, "12+20=", 2) | Out-Null
$cell = $t.Cell(10, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("97-1=", $true, $false, $false, $false, $false, $true, 1, $false, "86-79=", 2) | Out-Null
$cell = $t.Cell(10, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("0+16=", $true, $false, $false, $false, $false, $true, 1, $false, "60-54=", 2) | Out-Null
$cell = $t.Cell(10, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("68-64=", $true, $false, $false, $false, $false, $true, 1, $false, "71+3=", 2) | Out-Null
$cell = $t.Cell(10, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("89-53=", $true, $false, $false, $false, $false, $true, 1, $false, "96-47=", 2) | Out-Null
$cell = $t.Cell(10, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("57-36=", $true, $false, $false, $false, $false, $true, 1, $false, "23-14=", 2) | Out-Null
$cell = $t.Cell(11, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("11-1=", $true, $false, $false, $false, $false, $true, 1, $false, "82-54=", 2) | Out-Null
$cell = $t.Cell(11, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("98-76=", $true, $false, $false, $false, $false, $true, 1, $false, "13+19=", 2) | Out-Null
$cell = $t.Cell(11, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("57-16=", $true, $false, $false, $false, $false, $true, 1, $false, "45+54=", 2) | Out-Null
$cell = $t.Cell(11, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("99-63=", $true, $false, $false, $false, $false, $true, 1, $false, "23-21=", 2) | Out-Null
$cell = $t.Cell(11, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("81-73=", $true, $false, $false, $false, $false, $true, 1, $false, "57-18=", 2) | Out-Null
$cell = $t.Cell(12, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("79-52=", $true, $false, $false, $false, $false, $true, 1, $false, "25+27=", 2) | Out-Null
$cell = $t.Cell(12, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("3+80=", $true, $false, $false, $false, $false, $true, 1, $false, "81+1=", 2) | Out-Null
$cell = $t.Cell(12, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("2+78=", $true, $false, $false, $false, $false, $true, 1, $false, "35+12=", 2) | Out-Null
$cell = $t.Cell(12, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("67+16=", $true, $false, $false, $false, $false, $true, 1, $false, "47+17=", 2) | Out-Null
$cell = $t.Cell(12, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("52-44=", $true, $false, $false, $false, $false, $true, 1, $false, "48+4=", 2) | Out-Null
$cell = $t.Cell(13, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("11+43=", $true, $false, $false, $false, $false, $true, 1, $false, "63-24=", 2) | Out-Null
$cell = $t.Cell(13, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("32+35=", $true, $false, $false, $false, $false, $true, 1, $false, "98-7=", 2) | Out-Null
$cell = $t.Cell(13, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("33+14=", $true, $false, $false, $false, $false, $true, 1, $false, "96-3=", 2) | Out-Null
$cell = $t.Cell(13, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("65-39=", $true, $false, $false, $false, $false, $true, 1, $false, "51+11=", 2) | Out-Null
$cell = $t.Cell(13, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("22+74=", $true, $false, $false, $false, $false, $true, 1, $false, "11+60=", 2) | Out-Null
$cell = $t.Cell(14, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("52+19=", $true, $false, $false, $false, $false, $true, 1, $false, "64-22=", 2) | Out-Null
$cell = $t.Cell(14, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("25+57=", $true, $false, $false, $false, $false, $true, 1, $false, "86-73=", 2) | Out-Null
$cell = $t.Cell(14, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("34-3=", $true, $false, $false, $false, $false, $true, 1, $false, "82-11=", 2) | Out-Null
$cell = $t.Cell(14, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("80-41=", $true, $false, $false, $false, $false, $true, 1, $false, "38+31=", 2) | Out-Null
$cell = $t.Cell(14, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("47-29=", $true, $false, $false, $false, $false, $true, 1, $false, "76-46=", 2) | Out-Null
$cell = $t.Cell(15, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("73-67=", $true, $false, $false, $false, $false, $true, 1, $false, "32-10=", 2) | Out-Null
$cell = $t.Cell(15, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("98-4=", $true, $false, $false, $false, $false, $true, 1, $false, "26-21=", 2) | Out-Null
$cell = $t.Cell(15, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("12+62=", $true, $false, $false, $false, $false, $true, 1, $false, "67-12=", 2) | Out-Null
$cell = $t.Cell(15, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("45-34=", $true, $false, $false, $false, $false, $true, 1, $false, "48+41=", 2) | Out-Null
$cell = $t.Cell(15, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("25+72=", $true, $false, $false, $false, $false, $true, 1, $false, "54-41=", 2) | Out-Null
$cell = $t.Cell(16, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("70-57=", $true, $false, $false, $false, $false, $true, 1, $false, "85-19=", 2) | Out-Null
$cell = $t.Cell(16, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("81-74=", $true, $false, $false, $false, $false, $true, 1, $false, "38-34=", 2) | Out-Null
$cell = $t.Cell(16, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("37+55=", $true, $false, $false, $false, $false, $true, 1, $false, "15+34=", 2) | Out-Null
$cell = $t.Cell(16, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("1+26=", $true, $false, $false, $false, $false, $true, 1, $false, "66-38=", 2) | Out-Null
$cell = $t.Cell(16, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("6+39=", $true, $false, $false, $false, $false, $true, 1, $false, "64-1=", 2) | Out-Null
$cell = $t.Cell(17, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("44-27=", $true, $false, $false, $false, $false, $true, 1, $false, "64+4=", 2) | Out-Null
$cell = $t.Cell(17, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("52-29=", $true, $false, $false, $false, $false, $true, 1, $false, "3-0=", 2) | Out-Null
$cell = $t.Cell(17, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("39+37=", $true, $false, $false, $false, $false, $true, 1, $false, "17+13=", 2) | Out-Null
$cell = $t.Cell(17, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("65-63=", $true, $false, $false, $false, $false, $true, 1, $false, "24-10=", 2) | Out-Null
$cell = $t.Cell(17, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("93+6=", $true, $false, $false, $false, $false, $true, 1, $false, "16-0=", 2) | Out-Null
$cell = $t.Cell(18, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("90-81=", $true, $false, $false, $false, $false, $true, 1, $false, "95-27=", 2) | Out-Null
$cell = $t.Cell(18, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("30-0=", $true, $false, $false, $false, $false, $true, 1, $false, "18-6=", 2) | Out-Null
$cell = $t.Cell(18, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("62+7=", $true, $false, $false, $false, $false, $true, 1, $false, "74+0=", 2) | Out-Null
$cell = $t.Cell(18, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("16+42=", $true, $false, $false, $false, $false, $true, 1, $false, "10+35=", 2) | Out-Null
$cell = $t.Cell(18, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("89-46=", $true, $false, $false, $false, $false, $true, 1, $false, "81-6=", 2) | Out-Null
$cell = $t.Cell(19, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("32-24=", $true, $false, $false, $false, $false, $true, 1, $false, "40+0=", 2) | Out-Null
$cell = $t.Cell(19, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("31+10=", $true, $false, $false, $false, $false, $true, 1, $false, "92-40=", 2) | Out-Null
$cell = $t.Cell(19, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("9+35=", $true, $false, $false, $false, $false, $true, 1, $false, "70+0=", 2) | Out-Null
$cell = $t.Cell(19, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("22+52=", $true, $false, $false, $false, $false, $true, 1, $false, "12+84=", 2) | Out-Null
$cell = $t.Cell(19, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("46+30=", $true, $false, $false, $false, $false, $true, 1, $false, "56+2=", 2) | Out-Null
$cell = $t.Cell(20, 1)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("91-15=", $true, $false, $false, $false, $false, $true, 1, $false, "0+54=", 2) | Out-Null
$cell = $t.Cell(20, 2)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("21+5=", $true, $false, $false, $false, $false, $true, 1, $false, "90-19=", 2) | Out-Null
$cell = $t.Cell(20, 3)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("79+1=", $true, $false, $false, $false, $false, $true, 1, $false, "43-10=", 2) | Out-Null
$cell = $t.Cell(20, 4)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("33+32=", $true, $false, $false, $false, $false, $true, 1, $false, "81-10=", 2) | Out-Null
$cell = $t.Cell(20, 5)
$r0 = $cell.Range
$rr = $d.Range($r0.Start, $r0.End - 1)
$rr.Find.Execute("81-4=", $true, $false, $false, $false, $false, $true, 1, $false, "32+2=", 2) | Out-Null
